$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 268, pushing the existing rows 268-308 down to 269-309.
$ws.Rows.Item(268).Insert()

# Populate the newly inserted row 268 with a new price record. This mirrors the
# previous row 268 contents (same product/region/grade/price-per-unit etc.) but
# with an updated report date (D) and a new reported volume (M).
$ws.Range("A268").Value = 10
$ws.Range("B268").Value = "Vega Modelo de Temuco"
$ws.Range("C268").Value = "La Araucanía"
$ws.Range("D268").Value = 45131
$ws.Range("E268").Value = 9
$ws.Range("F268").Value = "Fruta"
$ws.Range("G268").Value = 100104
$ws.Range("H268").Value = "Frutos de pepita"
$ws.Range("I268").Value = 100104003
$ws.Range("J268").Value = "Membrillo"
$ws.Range("K268").Value = "Champion"
$ws.Range("L268").Value = "Primera"
$ws.Range("M268").Value = 100
$ws.Range("N268").Value = 16000
$ws.Range("O268").Value = 16000
$ws.Range("P268").Value = 16000
$ws.Range("Q268").Value = "$/bandeja 18 kilos granel"
$ws.Range("R268").Value = "Región de O'Higgins"
$ws.Range("S268").Value = 889
$ws.Range("T268").Value = 18
